$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column M (idContratoIXC), shifting M->N, N->O, O->P
$ws.Range("M1").EntireColumn.Insert()

# New header cell M1 = nomePlano, copy style from neighboring header cell (L1)
$ws.Range("M1").Value = "nomePlano"
$ws.Range("M1").Style = $ws.Range("L1").Style

# New data cell M2 = nomePlano value
$ws.Range("M2").Value = "100MB_SPEEDY"

# Update idContratoIXC value (now in column N) and B2 client id
$ws.Range("N2").Value = "154047"
$ws.Range("B2").Value = "117696"

# Update logRetorno JSON text (now in column P)
$ws.Range("P2").Value = "{'type': 'success', 'message': 'Registro inserido com sucesso!', 'id': '154047', 'atualiza_campos': [{'tipo': 'i', 'campo': 'data_desistencia', 'valor': ''}, {'tipo': 's', 'campo': 'status', 'valor': 'P'}, {'tipo': 'i', 'campo': 'id_cliente', 'valor': '117696'}, {'tipo': 'i', 'campo': 'data_ativacao', 'valor': ''}, {'tipo': 'd', 'campo': 'data_renovacao', 'valor': ''}, {'tipo': 'd', 'campo': 'nao_avisar_ate', 'valor': ''}, {'tipo': 'd', 'campo': 'nao_bloquear_ate', 'valor': ''}, {'tipo': 'd', 'campo': 'data_cancelamento', 'valor': ''}, {'tipo': 'd', 'campo': 'dt_ult_bloq_manual', 'valor': ''}, {'tipo': 'd', 'campo': 'data_cadastro_sistema', 'valor': '2025-03-19'}]}"
